# Add "Errors" and "Warnings" sheets with diagnostic messages to the
# "Classes" workbook, and move the active tab selection to "Warnings".

$wb = $excel.ActiveWorkbook
$classes = $wb.Worksheets.Item("Classes")

# --- Add the "Errors" sheet (placed right after "Classes") -----------------
$errors = $wb.Worksheets.Add($null, $classes)
$errors.Name = "Errors"

# The source strings start and end with a literal apostrophe. A leading
# apostrophe typed/assigned straight into Value is treated as Excel's
# "quote prefix" and is not stored as text, so build each string via a
# scratch formula (CHAR(39)/CHAR(34) for the quote characters) and paste
# the computed value back in as a literal so the apostrophe survives.
$scratch = $errors.Range("Z1")

$scratch.Formula = '=CHAR(39)&"Sheet "&CHAR(34)&"Classes"&CHAR(34)&" Row: 1 Column "&CHAR(34)&"A"&CHAR(34)&" in the header is not labeled as "&CHAR(34)&"DDBNNN"&CHAR(34)&CHAR(39)&","'
$scratch.Copy()
$errors.Range("A1").PasteSpecial(-4163)

$scratch.Formula = '=CHAR(39)&"Sheet "&CHAR(34)&"Classes"&CHAR(34)&" Row: 1 Column "&CHAR(34)&"B"&CHAR(34)&" in the header is not labeled as "&CHAR(34)&"TITLE"&CHAR(34)&CHAR(39)&","'
$scratch.Copy()
$errors.Range("A2").PasteSpecial(-4163)

$scratch.Formula = '=CHAR(39)&"Sheet "&CHAR(34)&"Classes"&CHAR(34)&" Row: 1 Column "&CHAR(34)&"C"&CHAR(34)&" in the header is not labeled as "&CHAR(34)&"OFF CLS"&CHAR(34)&CHAR(39)&","'
$scratch.Copy()
$errors.Range("A3").PasteSpecial(-4163)

$scratch.Formula = '=CHAR(39)&"Sheet "&CHAR(34)&"Classes"&CHAR(34)&" Row: 1 Column "&CHAR(34)&"D"&CHAR(34)&" in the header is not labeled as "&CHAR(34)&"SUB CLASSES"&CHAR(34)&CHAR(39)&","'
$scratch.Copy()
$errors.Range("A4").PasteSpecial(-4163)

$scratch.ClearContents()

$errors.Columns.Item(1).ColumnWidth = 62
$errors.Range("A11").Select()

# --- Add the "Warnings" sheet (placed right after "Errors") ----------------
$warnings = $wb.Worksheets.Add($null, $errors)
$warnings.Name = "Warnings"
$warnings.Range("D43").Select()

# --- "Classes" is no longer the selected tab --------------------------------
$classes.Select()

# --- "Warnings" becomes the active sheet/tab --------------------------------
$warnings.Activate()
$warnings.Select()
